$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

# Update PLC live data values as of 2025-10-13 13:46:17
$ws.Range("C2").Value = 260
$ws.Range("C3").Value = 158715
$ws.Range("C4").Value = 149771
$ws.Range("C8").Value = 64
